$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Q8" header in J1: copy formatting from I1 (bold/border/center style), then set value
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").Value = "Q8"

# Re-simulated rt_data / recomputed naive errors for each quarter row
# Row 2
$ws.Range("B2").Value = 0.5069022353042435
$ws.Range("C2").Value = -1.768726087092867
$ws.Range("D2").Value = 0.5753591407540363
$ws.Range("E2").Value = -0.04800876514786478
$ws.Range("F2").Value = 5.462709330733921
$ws.Range("G2").Value = 0.497888772914498
$ws.Range("H2").Value = 0.3521538662529956

# Row 3
$ws.Range("B3").Value = -1.813005556311658
$ws.Range("C3").Value = 0.5310796715352453
$ws.Range("D3").Value = -0.09228823436665579
$ws.Range("E3").Value = 5.418429861515129
$ws.Range("F3").Value = 0.453609303695707
$ws.Range("G3").Value = 0.3078743970342046

# Row 4
$ws.Range("B4").Value = 0.8445005399295014
$ws.Range("C4").Value = 0.2211326340276003
$ws.Range("D4").Value = 5.731850729909385
$ws.Range("E4").Value = 0.767030172089963
$ws.Range("F4").Value = 0.6212952654284607
$ws.Range("G4").Value = 0.4971820485470204
$ws.Range("H4").Value = 1.271373231878403
$ws.Range("I4").Value = -0.190180343825102
$ws.Range("J4").Value = 0.410128014204804

# Row 5
$ws.Range("B5").Value = -0.3768185295679458
$ws.Range("C5").Value = 5.133899566313839
$ws.Range("D5").Value = 0.169079008494417
$ws.Range("E5").Value = 0.0233441018329146
$ws.Range("F5").Value = -0.1007691150485257
$ws.Range("G5").Value = 0.6734220682828564
$ws.Range("H5").Value = -0.7881315074206481
$ws.Range("I5").Value = -0.187823149390742

# Row 6
$ws.Range("B6").Value = 5.029545179712666
$ws.Range("C6").Value = 0.06472462189324377
$ws.Range("D6").Value = -0.08101028476825861
$ws.Range("E6").Value = -0.2051235016496989
$ws.Range("F6").Value = 0.5690676816816833
$ws.Range("G6").Value = -0.8924858940218212
$ws.Range("H6").Value = -0.2921775359919152

# Row 7
$ws.Range("B7").Value = -0.137413759447277
$ws.Range("C7").Value = -0.2831486661087794
$ws.Range("D7").Value = -0.4072618829902197
$ws.Range("E7").Value = 0.3669293003411625
$ws.Range("F7").Value = -1.094624275362342
$ws.Range("G7").Value = -0.494315917332436

# Row 8
$ws.Range("B8").Value = -0.1897021710626319
$ws.Range("C8").Value = -0.3138153879440722
$ws.Range("D8").Value = 0.46037579538731
$ws.Range("E8").Value = -1.001177780316195
$ws.Range("F8").Value = -0.4008694222862885
$ws.Range("G8").Value = -0.3825116209597155
$ws.Range("H8").Value = 0.3758429424091532
$ws.Range("I8").Value = -0.7584425099968086

# Row 9
$ws.Range("B9").Value = -0.7538669319268549
$ws.Range("C9").Value = 0.0203242514045272
$ws.Range("D9").Value = -1.441229324298977
$ws.Range("E9").Value = -0.8409209662690713
$ws.Range("F9").Value = -0.8225631649424983
$ws.Range("G9").Value = -0.06420860157362956
$ws.Range("H9").Value = -1.198494053979591

# Row 10
$ws.Range("B10").Value = 0.4503073135717173
$ws.Range("C10").Value = -1.011246262131787
$ws.Range("D10").Value = -0.4109379041018812
$ws.Range("E10").Value = -0.3925801027753081
$ws.Range("F10").Value = 0.3657744605935606
$ws.Range("G10").Value = -0.7685109918124013

# Row 11
$ws.Range("B11").Value = -0.9542016326501267
$ws.Range("C11").Value = -0.3538932746202207
$ws.Range("D11").Value = -0.3355354732936476
$ws.Range("E11").Value = 0.4228190900752211
$ws.Range("F11").Value = -0.7114663623307408

# Row 12
$ws.Range("B12").Value = -0.3135312646359353
$ws.Range("C12").Value = -0.2951734633093623
$ws.Range("D12").Value = 0.4631811000595064
$ws.Range("E12").Value = -0.6711043523464554

# Row 13
$ws.Range("B13").Value = -0.2451295611021919
$ws.Range("C13").Value = 0.5132250022666769
$ws.Range("D13").Value = -0.621060450139285

# Row 14
$ws.Range("B14").Value = 0.8266790722624406
$ws.Range("C14").Value = -0.3076063801435212

# Row 15
$ws.Range("B15").Value = -0.3469781724577359

Write-Host "done"
